# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the old Strike# values in column G, rows 2-25
$newValues = @{
    2  = 2
    3  = 5
    4  = 4
    5  = 8
    6  = 9
    7  = 7
    8  = 5
    9  = 10
    10 = 12
    11 = 3
    12 = 6
    13 = 6
    14 = 10
    15 = 9
    16 = 5
    17 = 7
    18 = 13
    19 = 10
    20 = 11
    21 = 3
    22 = 6
    23 = 8
    24 = 10
    25 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
